$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Select the entire second row (the duplicate data row) and delete it,
# shifting all subsequent rows up by one.
$row = $ws.Rows.Item(2)
$row.Select()
$row.Delete()

# Update selection to mirror post-delete state (entire row 2 selected).
$ws.Range("A2:XFD2").Select()
